$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -14.0048
$ws.Range("C6").Value = -11.65900000000001
$ws.Range("C7").Value = -12.06460000000001
$ws.Range("C8").Value = -11.61859999999999
$ws.Range("C16").Value = -11.6518
$ws.Range("C20").Value = -14.5997
$ws.Range("C21").Value = -13.22900000000001
